# PyramidRecruiterDashboard.xlsx - "added few more QA tasks"
#
# Inserts three new QA task rows (Bussiness Requirement Creation, Test Plan,
# Automation Framework set-up) right above the existing "Test Case Execution"
# row, pushing the remainder of the task table down by three rows. Also
# updates the active view to show the newly inserted area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 blank rows before row 39 (old row 39 -> new row 42, etc.)
$ws.Rows("39:41").Insert()

# The freshly inserted rows pick up a plain default style; copy the
# formatting (borders / wrap-text) from the row just below them (the
# row that used to be row 39, now shifted to row 42) so they match the
# rest of the table.
$ws.Range("A42:H42").Copy()
$ws.Range("A39:H41").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New row 39 - Bussiness Requirement Creation
$ws.Range("C39").Value = "Bussiness Requirement Creation"
$ws.Range("E39").Value = "Vijay"
$ws.Range("G39").Value = "Done"

# New row 40 - Test Plan
$ws.Range("C40").Value = "Test Plan"
$ws.Range("E40").Value = "Vijay"
$ws.Range("G40").Value = "Done"

# New row 41 - Automation Framework set-up
$ws.Range("C41").Value = "Automation Framework set-up"
$ws.Range("E41").Value = "Vijay"
$ws.Range("G41").Value = "Done"

# Update the active view: scrolled down to the new task rows with
# C42 ("Test Case Execution") selected.
$ws.Activate()
$ws.Range("C42").Select()
